$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block: rows 2-6 (trial rows), columns A-J
# A=trialTrain, B=x_fixStart, C=y_fixStart, D=x_corrSteps, E=y_corrSteps,
# F=x_nrSteps, G=y_nrSteps, H=alienID, I=praclen, J=version
$data = @(
    @(1, 0, 8, 2, 4, 2, -4, 45, 5, "train_dim2_1"),
    @(2, 1, 7, 2, 2, 1, -5, 56, 5, "train_dim2_1"),
    @(3, 1, 9, 6, 8, 5, -1, 12, 5, "train_dim2_1"),
    @(4, 0, 6, 3, 3, 3, -3, 34, 5, "train_dim2_1"),
    @(5, 3, 9, 7, 7, 4, -2, 23, 5, "train_dim2_1")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

# Update selection to match final state
$ws.Range("I1").Select()
